# BOT; UPDATE DATA
# Adds one new daily data row (2020-05-15 / serial 43966) to the
# "相談件数" sheet, pushing the trailing footnote row down by one, and
# updates the sheet/workbook view state (dimension, print area, frozen
# pane scroll position, active selection) to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Insert a new row 111 (shifts the old row 111 -> row 112) -----------
# Inserting via a full row range copies down the formatting/number
# styles of the row above, which already match what the new data row
# needs (date / count / count / right-aligned count columns).
[void]($ws.Range("A111:E111").Insert())

$ws.Range("A111").Value = 43966
$ws.Range("B111").Value = 216
$ws.Range("C111").Value = 37290
$ws.Range("D111").Value = 54
$ws.Range("E111").Value = 7584

# --- Print area -----------------------------------------------------------
# Source workbook keeps the print area one row taller than the actual
# data range; follow the same convention (new last data row is 112).
$ws.PageSetup.PrintArea = '$A$1:$E$113'

# --- Frozen pane / selection state -----------------------------------------
$win = $excel.ActiveWindow

# Re-establish the freeze (col A / row 1) and then scroll the
# bottom-right pane so row 95 is the first visible row, matching the
# new topLeftCell of the frozen pane.
[void]($win.FreezePanes = $false)
[void]($ws.Range("B2").Select())
[void]($win.FreezePanes = $true)
[void]($excel.Goto($ws.Range("C95"), $true))

# Final active cell / selection, as left by the author after appending
# the new row.
[void]($ws.Range("E109").Select())
